$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15816
$ws.Range("D2").Value = 0.4392357253943568

$ws.Range("C3").Value = 20192
$ws.Range("D3").Value = 0.5607642746056432

$ws.Range("C4").Value = 24584
$ws.Range("D4").Value = 0.6827371695178849

$ws.Range("C5").Value = 6941
$ws.Range("D5").Value = 0.1927627193956898

$ws.Range("C6").Value = 4483
$ws.Range("D6").Value = 0.1245001110864252

$ws.Range("C7").Value = 10706
$ws.Range("D7").Value = 0.2973228171517441

$ws.Range("C8").Value = 18426
$ws.Range("D8").Value = 0.5117196178626972

$ws.Range("C9").Value = 6876
$ws.Range("D9").Value = 0.1909575649855588
